$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2..4) down to (3..5), working bottom-up so we
# don't overwrite data we still need to read.
for ($r = 4; $r -ge 2; $r--) {
    $dest = $r + 1
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($dest, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
    $ws.Cells.Item($dest, 4).NumberFormat = $ws.Cells.Item($r, 4).NumberFormat
}

# Populate the new row 2 with the newly reported weekly price entry.
$ws.Cells.Item(2, 1).Value = 10
$ws.Cells.Item(2, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(2, 3).Value = "La Araucanía"
$ws.Cells.Item(2, 4).Value = 44453
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
$ws.Cells.Item(2, 5).Value = 9
$ws.Cells.Item(2, 6).Value = 100112042
$ws.Cells.Item(2, 7).Value = "Locoto"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 20
$ws.Cells.Item(2, 11).Value = 2300
$ws.Cells.Item(2, 12).Value = 2300
$ws.Cells.Item(2, 13).Value = 2300
$ws.Cells.Item(2, 14).Value = "$/kilo"
$ws.Cells.Item(2, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(2, 16).Value = 2300
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"
